$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 8-18 (11 species-observation rows) get their row contents cyclically
# rotated: the data that used to live in the last two rows (17,18) moves up
# to become rows 8,9, and every other row's data shifts down by two rows
# (old 8->new 10, old 9->new 11, ... old 16->new 18). Only the columns whose
# values actually differ from row to row are touched (A,B,D,E,F,G,H,I,J,K,L
# and Q,R) so untouched columns (dates, location text, etc. -- identical on
# every one of these rows) are left completely alone.

$firstRow = 8
$lastRow = 18
$rowCount = $lastRow - $firstRow + 1

function Get-RotatedArray($range) {
    $data = $range.Value2
    $colCount = $data.GetLength(1)
    $newData = New-Object 'object[,]' $rowCount, $colCount
    for ($i = 1; $i -le $rowCount; $i++) {
        $srcIndex = (($i - 1 - 2) % $rowCount + $rowCount) % $rowCount + 1
        for ($j = 1; $j -le $colCount; $j++) {
            $newData[$i - 1, $j - 1] = $data[$srcIndex, $j]
        }
    }
    return $newData
}

$rangeAL = $ws.Range("A${firstRow}:L${lastRow}")
$newAL = Get-RotatedArray $rangeAL
$rangeAL.Value = $newAL

$rangeQR = $ws.Range("Q${firstRow}:R${lastRow}")
$newQR = Get-RotatedArray $rangeQR
$rangeQR.Value = $newQR
